$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Connectivity")
$ws1.Range("A1").Value = [double]"0.005857887150227237"
$ws1.Range("B1").Value = [double]"0.005033460283055886"
$ws1.Range("C1").Value = [double]"0.005339797749090689"
$ws1.Range("D1").Value = [double]"0.0053896520392816985"
$ws1.Range("E1").Value = [double]"2.358906964794187E-4"
$ws1.Range("F1").Value = [double]"3.1679263854308676E307"
$ws1.Range("G1").Value = [double]"1.3589459247596537E308"
$ws1.Range("H1").Value = [double]"7.482971686277037E307"
$ws1.Range("I1").Value = [double]"1.2483570308335637E308"
$ws1.Range("J1").Value = [double]"1.5131464316013779E308"

$ws1.Range("A2").Value = 27.0
$ws1.Range("B2").Value = 0.0
$ws1.Range("C2").Value = 6.0
$ws1.Range("D2").Value = 8.0

$ws2 = $wb.Worksheets.Item("Is Isolated")
$ws2.Range("A1").Value = [double]"2.8560776309376505E-4"
$ws2.Range("B1").Value = [double]"0.00372002482552455"
$ws2.Range("C1").Value = [double]"0.006483710142371962"
$ws2.Range("D1").Value = [double]"0.003743810130259544"
$ws2.Range("E1").Value = [double]"0.004176234691075986"
$ws2.Range("F1").Value = [double]"1.2209648322088105E308"
$ws2.Range("G1").Value = [double]"1.3661546677709258E307"
$ws2.Range("H1").Value = [double]"3.461662078480697E307"
$ws2.Range("I1").Value = [double]"1.6768333882002263E308"
$ws2.Range("J1").Value = [double]"1.9304363359560103E307"

$ws2.Range("C2").Value = 394.0
$ws2.Range("E2").Value = 500.0

$ws3 = $wb.Worksheets.Item("Diameter")
$ws3.Range("A1").Value = [double]"0.01553666732814921"
$ws3.Range("B1").Value = [double]"0.04467448688954686"
$ws3.Range("C1").Value = [double]"0.07668788255007811"
$ws3.Range("D1").Value = [double]"0.10937757571686146"
$ws3.Range("E1").Value = [double]"0.054827703426204036"
$ws3.Range("F1").Value = [double]"9.671823698659732E307"
$ws3.Range("G1").Value = [double]"9.941426995253211E306"
$ws3.Range("H1").Value = [double]"1.6920586194097304E308"
$ws3.Range("I1").Value = [double]"1.5914341940751387E308"
$ws3.Range("J1").Value = [double]"8.300351259686377E307"

$ws3.Range("A2").Value = 20.0
$ws3.Range("B2").Value = 15.0
$ws3.Range("C2").Value = 14.0
$ws3.Range("D2").Value = 10.0
$ws3.Range("E2").Value = 15.0
$ws3.Range("F2").Value = 5.0
$ws3.Range("G2").Value = 5.0
$ws3.Range("H2").Value = 5.0
$ws3.Range("I2").Value = 5.0
$ws3.Range("J2").Value = 5.0
